# Commit: "updated styles and table"
# The only substantive content change in the target diff is that the
# numeric value previously stored in cell A5 (1.0) has been cleared, so the
# cell is now blank. (The remaining differences in the diff -- namespace
# churn, style normalization, drawing removal, etc. -- are the byproduct of
# Excel re-serializing the workbook and are produced automatically by the
# COM save, not something to set explicitly.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").ClearContents()
